# SE_Est.xlsx : first round of results with SMM.
# Expand the small 5-col x 6-row summary table into a 12-col x 5-row table:
#   - new header columns for rho / sigma SE's (besides the existing SPF / SCE ones)
#   - new row labels (FEVar / FEATV / DisgATV / DisgVar) replacing the old
#     (Forecast / FE / Disg / Var) labels
#   - new numeric data matching the new layout

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header row (row 1): A:D are the plain numbers 0,1,2,3 (style 1);
# E:L are shared-string labels (style 1) built from the existing label
# text reused across pairs of columns.
# ---------------------------------------------------------------------
$ws.Range("A1").Value = 0
$ws.Range("B1").Value = 1
$ws.Range("C1").Value = 2
$ws.Range("D1").Value = 3

$ws.Range("E1").Value = 'SE: $\hat\lambda_{SPF}$(Q)'
$ws.Range("F1").Value = 'SE: $\hat\lambda_{SPF}$(Q)'
$ws.Range("G1").Value = 'SE: $\rho$'
$ws.Range("H1").Value = 'SE: $\sigma$'
$ws.Range("I1").Value = 'SE: $\hat\lambda_{SCE}$(M)'
$ws.Range("J1").Value = 'SE: $\hat\lambda_{SCE}$(M)'
$ws.Range("K1").Value = 'SE: $\rho$'
$ws.Range("L1").Value = 'SE: $\sigma$'

# D1 used to be a shared string ("SPF" label) - it is now a plain number,
# so make sure any leftover string-type flag is gone by re-setting it.
$ws.Range("D1").Value = 3

# F1:L1 are brand-new cells outside the old A1:E6 used range, so they do
# not automatically inherit the header style (bold + border + centered).
# Copy the formatting from the existing styled header cell A1 so the
# workbook's styles.xml is reused as-is (no new cellXfs entries).
$ws.Range("A1").Copy() | Out-Null
$ws.Range("F1:L1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Row 2
# ---------------------------------------------------------------------
$ws.Range("A2").Value = "FEVar"
$ws.Range("B2").Value = "FEATV"
$ws.Range("C2").Clear()
$ws.Range("D2").Clear()
$ws.Range("E2").Value = 0.47
$ws.Range("F2").Value = 0.36
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 0.08
$ws.Range("I2").Value = 0.2
$ws.Range("J2").Value = 0.59
$ws.Range("K2").Value = 0.99
$ws.Range("L2").Value = 0.08

# ---------------------------------------------------------------------
# Row 3
# ---------------------------------------------------------------------
$ws.Range("A3").Value = "FEVar"
$ws.Range("B3").Value = "DisgATV"
$ws.Range("C3").Value = "DisgVar"
$ws.Range("D3").Clear()
$ws.Range("E3").Value = 0.27
$ws.Range("F3").Value = 0.38
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 0.11
$ws.Range("I3").Value = 0.2
$ws.Range("J3").Value = 0.5600000000000001
$ws.Range("K3").Value = 0.98
$ws.Range("L3").Value = 0.08

# ---------------------------------------------------------------------
# Row 4
# ---------------------------------------------------------------------
$ws.Range("A4").Value = "FEVar"
$ws.Range("B4").Value = "FEATV"
$ws.Range("C4").Value = "DisgVar"
$ws.Range("D4").Value = "DisgATV"
$ws.Range("E4").Value = 0.47
$ws.Range("F4").Value = 0.36
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 0.1
$ws.Range("I4").Value = 0.2
$ws.Range("J4").Value = 0.59
$ws.Range("K4").Value = 0.99
$ws.Range("L4").Value = 0.08

# ---------------------------------------------------------------------
# Row 5 (replaces the old row 5 + row 6)
# ---------------------------------------------------------------------
$ws.Range("A5").Value = "FEVar"
$ws.Range("B5").Value = "FEATV"
$ws.Range("C5").Value = "DisgVar"
$ws.Range("D5").Value = "DisgATV"
$ws.Range("E5").Value = 0.47
$ws.Range("F5").Value = 0.36
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = 0.1
$ws.Range("I5").Value = 0.2
$ws.Range("J5").Value = 0.59
$ws.Range("K5").Value = 0.99
$ws.Range("L5").Value = 0.08

# ---------------------------------------------------------------------
# The old sheet had a 6th row - remove it entirely, the new table is
# only 5 rows tall (A1:L5).
# ---------------------------------------------------------------------
$ws.Range("A6:L6").Clear()
